$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147 (existing rows 147-161 shift down to 148-162)
$ws.Rows(147).Insert()

# Populate the newly inserted row 147 with the new weekly record
$ws.Range("A147").Value = 10
$ws.Range("B147").Value = "Vega Modelo de Temuco"
$ws.Range("C147").Value = "La Araucanía"
$ws.Range("D147").Value = 44504
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100102
$ws.Range("H147").Value = "Cítricos"
$ws.Range("I147").Value = 100102006
$ws.Range("J147").Value = "Pomelo"
$ws.Range("K147").Value = "Start Ruby"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 140
$ws.Range("N147").Value = 12000
$ws.Range("O147").Value = 12000
$ws.Range("P147").Value = 12000
$ws.Range("Q147").Value = "$/bandeja 15 kilos granel"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 800
$ws.Range("T147").Value = 15
